# ---------------------------------------------------------------------------
# relatorio_meta_1.docx -- "MAIN: Add minor improvements to report."
#
# Changes applied:
#  1. Rename several internal style IDs (table styles + character styles)
#     to their PT-PT localized equivalents (display names / w:name stay the
#     same except for the two "...Char" companion character styles, whose
#     NameLocal also changes).
#  2. Split "Máscara de Rede" into three runs so the last word becomes
#     lower-case "rede" (x4 occurrences).
#  3. Italicise the word "Router" in five table-cell labels
#     ("Router R2 – ", "Router R1", "Router R2", "Router R3").
#  4. Add a caption ("Tabela 9 - Configuração do roteamento dos três
#     routers") to the previously empty paragraph that follows the last
#     table.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: pull the full OOXML package (all parts, incl. styles.xml) through
# a tiny range's WordOpenXML, edit it as plain text, then feed it back with
# InsertXML over the whole document. InsertXML merges in any style that is
# new to the package but does not strip styles that already existed, so any
# style whose *id* changed leaves behind a stale copy under its old id --
# those are removed explicitly afterwards via Styles.Item(old).Delete().
# ---------------------------------------------------------------------------

$probe = $d.Range(0, 1)
$xml = $probe.WordOpenXML

# --- style id renames (old id -> new id) ------------------------------------------------
$styleIdRenames = @(
    @("Tabelacomgrade", "TabelacomGrelha"),
    @("Fontepargpadro", "Tipodeletrapredefinidodopargrafo"),
    @("TabeladeLista3-nfase1", "TabeladeLista3-Destaque1"),
    @("TabeladeGrade4-nfase1", "TabeladeGrelha4-Destaque1"),
    @("TabeladeGrade2-nfase2", "TabeladeGrelha2-Destaque2"),
    @("TabeladeLista4-nfase2", "TabeladeLista4-Destaque2"),
    @("TabeladeLista4-nfase3", "TabeladeLista4-Destaque3"),
    @("TabeladeGrade4-nfase3", "TabeladeGrelha4-Destaque3"),
    @("CabealhoChar", "CabealhoCarter"),
    @("RodapChar", "RodapCarter")
)

foreach ($pair in $styleIdRenames) {
    $old = $pair[0]
    $new = $pair[1]
    $xml = $xml.Replace('w:styleId="' + $old + '"', 'w:styleId="' + $new + '"')
    $xml = $xml.Replace('w:val="' + $old + '"', 'w:val="' + $new + '"')
}

# --- the two localized display-name (w:name) updates that ride along with the id rename ---
$xml = $xml.Replace('<w:name w:val="Cabeçalho Char"/>', '<w:name w:val="Cabeçalho Caráter"/>')
$xml = $xml.Replace('<w:name w:val="Rodapé Char"/>', '<w:name w:val="Rodapé Caráter"/>')

# --- "Máscara de Rede" -> "Máscara de " + "r" + "ede"  (x4) ---
$mascaraOld = '<w:r w:rsidRPr="0002622B"><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Máscara de Rede</w:t></w:r>'
$mascaraNew = '<w:r w:rsidRPr="0002622B"><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">Máscara de </w:t></w:r><w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>r</w:t></w:r><w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>ede</w:t></w:r>'
$xml = $xml.Replace($mascaraOld, $mascaraNew)

# --- "Router R2 – " -> italic "Router" + " R2 – " ---
$r2dashOld = '<w:r w:rsidRPr="0002622B"><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">Router R2 – </w:t></w:r>'
$r2dashNew = '<w:r w:rsidRPr="0002622B"><w:rPr><w:i/><w:iCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Router</w:t></w:r><w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> R2 – </w:t></w:r>'
$xml = $xml.Replace($r2dashOld, $r2dashNew)

# --- "Router R1"/"Router R2"/"Router R3" -> italic "Router" + " R<n>" ---
foreach ($n in @("1", "2", "3")) {
    $old = '<w:r w:rsidRPr="006B2E47"><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Router R' + $n + '</w:t></w:r>'
    $new = '<w:r w:rsidRPr="006B2E47"><w:rPr><w:i/><w:iCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Router</w:t></w:r><w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> R' + $n + '</w:t></w:r>'
    $xml = $xml.Replace($old, $new)
}

# --- add the "Tabela 9 - ..." caption to the final (empty) paragraph ---
$finalParaOld = '<w:p w14:paraId="37D78F0D" w14:textId="77777777" w:rsidR="006B2E47" w:rsidRPr="00540025" w:rsidRDefault="006B2E47" w:rsidP="006B2E47"><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr></w:p>'

$captionRuns = ''
$captionRuns += '<w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">  </w:t></w:r>'
$captionRuns += '<w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="18"/><w:szCs w:val="18"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">Tabela </w:t></w:r>'
$captionRuns += '<w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="18"/><w:szCs w:val="18"/><w:u w:val="single"/></w:rPr><w:t>9</w:t></w:r>'
$captionRuns += '<w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>-</w:t></w:r>'
$captionRuns += '<w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> Configuração </w:t></w:r>'
$captionRuns += '<w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>do roteamento</w:t></w:r>'
$captionRuns += '<w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>'
$captionRuns += '<w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>d</w:t></w:r>'
$captionRuns += '<w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">os </w:t></w:r>'
$captionRuns += '<w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">três </w:t></w:r>'
$captionRuns += '<w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>routers</w:t></w:r>'

$finalParaNew = '<w:p w14:paraId="37D78F0D" w14:textId="77777777" w:rsidR="006B2E47" w:rsidRPr="00540025" w:rsidRDefault="006B2E47" w:rsidP="006B2E47"><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr>' + $captionRuns + '</w:p>'

$xml = $xml.Replace($finalParaOld, $finalParaNew)

# ---------------------------------------------------------------------------
# Step 2: write the edited package back, then drop the now-orphaned old
# style ids that InsertXML's merge semantics left behind.
# ---------------------------------------------------------------------------

$full = $d.Range()
$null = $full.InsertXML($xml)

foreach ($pair in $styleIdRenames) {
    $old = $pair[0]
    try {
        $oldStyle = $d.Styles.Item($old)
        $oldStyle.Delete()
    } catch {
        # already gone / never existed under the old id -- nothing to do
    }
}

Write-Output "done"
